$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# Update the BAU Guaranteed Dispatch Percentage (column B, year 2015) from 0 to 1
# for: onshore wind (row 6), solar PV (row 7), solar thermal (row 8), biomass (row 9),
# geothermal (row 10), offshore wind (row 14), municipal solid waste (row 17).
# Downstream years (columns C:AK) are driven by formulas referencing column B (or the
# prior column), so they recalculate automatically -- except row 17, which previously
# held hard-coded 0 values with no formulas and now needs them added.
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("B17").Value = 1

$ws.Range("C17").Formula = "=B17"
$ws.Range("D17:AK17").Formula = "=C17"

# Restore cursor/selection state: BGDPbES sheet remembers B6 was selected, then focus
# returns to the About sheet with C45 selected (About remains the active tab).
$ws.Activate()
$ws.Range("B6").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("C45").Select()
